$wb = $excel.ActiveWorkbook

# --- 1. Duplicate the existing "Sheet1" (536 2022 schedule) ---------------
# The copy is placed BEFORE the original, becomes the new "4536" sheet
# (keeps the old content untouched), while the original sheet is renamed
# to "536" and its contents updated to the new (2023) schedule.
$orig = $wb.Worksheets.Item(1)
$orig.Copy($orig)

$wb.Worksheets.Item(1).Name = "4536"
$wb.Worksheets.Item(2).Name = "536"

$ws = $wb.Worksheets.Item("536")

# --- 2. Update the header/row-2 date ---------------------------------------
$ws.Range("C2").Value = 45163

# --- 3. Row 4 (was "Pick Project Teams") -> "Rust 1: Setup + Vocab" --------
$ws.Range("B4").Value = "Rust 1: Setup + Vocab"
$ws.Range("C4").Value = 45167

# --- 4. Row 5 (was "What is a Language?") -> "Regular Expressions" ---------
$ws.Range("B5").Value = "Regular Expressions"
$ws.Range("C5").Value = 45170
$ws.Range("D5").Value = "A1"

# --- 5. Row 7 (was "Scala 1: Setup+Vocab") -> "Context Free Grammars" ------
$ws.Range("B7").Value = "Context Free Grammars"
$ws.Range("C7").Value = 45174

# --- 6. Row 8 (was "Regular Expressions") -> "Rust 2: Syntax Trees" --------
$ws.Range("B8").Value = "Rust 2: Syntax Trees"
$ws.Range("C8").Value = 45177
$ws.Range("D8").Value = "A2"

# --- 7. Row 10 (was "Context Free Grammars") -> "Rust 3: Interpreters" -----
$ws.Range("B10").Value = "Rust 3: Interpreters"
$ws.Range("C10").Value = 45181

# --- 8. Row 11 (was "Scala 2: Fastparse") -> "Operational Semantics" -------
$ws.Range("B11").Value = "Operational Semantics"
$ws.Range("C11").Value = 45184
$ws.Range("D11").Value = "A3"

# --- 9. Row 13 (was "Scala 3: Syntax Trees") -> "Types" --------------------
$ws.Range("B13").Value = "Types"
$ws.Range("C13").Value = 45188
$ws.Range("D13").ClearContents()

# --- 10. Row 14 (was "Scala 4: Interpreters") is removed entirely ----------
$ws.Range("A14").Clear()
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()

# --- 11. New row 15: "..." reading link -------------------------------------
$ws.Range("A15").Value = 9
$ws.Range("B15").Value = "https://pubmed.ncbi.nlm.nih.gov/16204405/"
$ws.Range("C15").Value = 45195

# --- 12. Row 16 (was "User Studies 1") -> "User Studies 2" -----------------
$ws.Range("B16").Value = "User Studies 2"
$ws.Range("C16").Value = 45198
$ws.Range("D16").Value = "A4"

# --- 13. Row 17 (was "User Studies 2") is removed entirely -----------------
$ws.Range("A17").Clear()
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()
$ws.Range("E17").Clear()

# --- 14. New row 18: "Do the Study" -----------------------------------------
$ws.Range("A18").Value = 11
$ws.Range("B18").Value = "Do the Study"
$ws.Range("C18").Value = 45202
$ws.Range("E18").Value = "MC"

# --- 15. Row 19 (was "User Studies 3") -> "Critical Code Studies" ----------
$ws.Range("B19").Value = "Critical Code Studies"
$ws.Range("C19").Value = 45205
$ws.Range("D19").ClearContents()

# --- 16. Row 21 (was "Do the Study") -> "Visual Arts" ----------------------
$ws.Range("B21").Value = "Visual Arts"
$ws.Range("C21").Value = 45209

# --- 17. Row 22 (was "Reflection + Planning") -> "Interactive Fiction" -----
$ws.Range("B22").Value = "Interactive Fiction"
$ws.Range("C22").Value = 45212
$ws.Range("D22").Value = "A5"

# --- 18. Rows 24,25,27,28: clear A/B, blank out date but keep date style ---
$ws.Range("A24").Clear()
$ws.Range("B24").Clear()
$ws.Range("C24").ClearContents()

$ws.Range("A25").Clear()
$ws.Range("B25").Clear()
$ws.Range("C25").ClearContents()

$ws.Range("A27").Clear()
$ws.Range("B27").Clear()
$ws.Range("C27").ClearContents()

$ws.Range("A28").Clear()
$ws.Range("B28").Clear()
$ws.Range("D28").Clear()
$ws.Range("C28").ClearContents()

# --- 19. Row 30: keep the hyperlink, drop the schedule columns -------------
$ws.Range("A30").Clear()
$ws.Range("B30").Clear()
$ws.Range("E30").Clear()
$ws.Range("C30").ClearContents()

# --- 20. Row 32 ---------------------------------------------------------
$ws.Range("A32").Clear()
$ws.Range("B32").Clear()
$ws.Range("E32").Clear()
$ws.Range("C32").ClearContents()

# --- 21. Row 33 ---------------------------------------------------------
$ws.Range("A33").Clear()
$ws.Range("B33").Clear()
$ws.Range("E33").Clear()
$ws.Range("C33").ClearContents()

# --- 22. Row 35 ---------------------------------------------------------
$ws.Range("A35").Clear()
$ws.Range("B35").Clear()
$ws.Range("E35").Clear()
$ws.Range("C35").ClearContents()

# --- 23. Row 36: everything except L36 is dropped (C36 fully cleared) -----
$ws.Range("A36").Clear()
$ws.Range("B36").Clear()
$ws.Range("C36").Clear()
$ws.Range("E36").Clear()
$ws.Range("H36").Clear()

# --- 24. New row 37 carries the link that used to sit on row 36's H -------
$ws.Range("C37").Value = 45188
$ws.Range("C37").ClearContents()
$ws.Range("H37").Value = "https://link.springer.com/content/pdf/10.1007/s00146-006-0050-9.pdf?pdf=button"

# --- 25. Row 38 ---------------------------------------------------------
$ws.Range("A38").Clear()
$ws.Range("B38").Clear()
$ws.Range("D38").Clear()
$ws.Range("E38").Clear()
$ws.Range("C38").ClearContents()

# --- 26. Row 39 is removed entirely (its H-link moves to row 40) ----------
$ws.Range("A39").Clear()
$ws.Range("B39").Clear()
$ws.Range("C39").Clear()
$ws.Range("E39").Clear()
$ws.Range("H39").Clear()

# --- 27. New row 40 carries the link that used to sit on row 39's H -------
$ws.Range("C40").Value = 1
$ws.Range("C40").ClearContents()
$ws.Range("H40").Value = "https://dl.acm.org/doi/pdf/10.1145/3393914.3395885"

# --- 28. Row 41 ---------------------------------------------------------
$ws.Range("A41").Clear()
$ws.Range("B41").Clear()
$ws.Range("E41").Clear()
$ws.Range("C41").ClearContents()

# --- 29. Row 42 is removed entirely -----------------------------------------
$ws.Range("A42").Clear()
$ws.Range("B42").Clear()
$ws.Range("C42").Clear()
$ws.Range("E42").Clear()
$ws.Range("P42").Clear()

# --- 30. Row 43: drop schedule columns, change Part label to "TBA" --------
$ws.Range("A43").Clear()
$ws.Range("B43").Clear()
$ws.Range("D43").Clear()
$ws.Range("C43").ClearContents()
$ws.Range("P43").Value = "TBA"

# --- 31. New row 44: "Project Celebration" now lives in the Part column ---
$ws.Range("C44").Value = 1
$ws.Range("C44").ClearContents()
$ws.Range("P44").Value = "Project Celebration"

# --- 32. Hyperlink bookkeeping ----------------------------------------------
# H30/H32/H33/H35/H38/H41 keep their original text/links (unchanged);
# H37 and H40 are brand-new hyperlinks (see above) that need the actual
# Hyperlinks collection entries, not just display text.
$ws.Hyperlinks.Add($ws.Range("H37"), "https://link.springer.com/content/pdf/10.1007/s00146-006-0050-9.pdf?pdf=button") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H40"), "https://dl.acm.org/doi/pdf/10.1145/3393914.3395885") | Out-Null

# --- 33. Sheet view / selection tweaks --------------------------------------
$wb.Worksheets.Item("4536").Range("B5").Select()
$ws.Range("B15").Select()
$ws.Activate()
